$d = $word.ActiveDocument

# Locate the paragraph "For responses to Email1:" so we can insert the new
# sub-bullet right after it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -match "For responses to Email1:") {
        $target = $p
        break
    }
}

# Insert a new paragraph immediately after the target paragraph.
$target.Range.InsertParagraphAfter()
$newPara = $target.Next()

# Set its text to the new checklist item.
$newPara.Range.Text = "Check that the user has responded with their eye gaze tracker model. If they have not, send a response re-requesting the information"

# Match the nested bullet level (w:ilvl=1) used by the sibling sub-items
# under this bullet (e.g. "Update spreadsheet with new participant ID").
$newPara.Range.ListFormat.ListLevelNumber = 2
